$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.344.53'
$ws.Range("E2").Value = '  -0.05%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.687.67'
$ws.Range("E3").Value = '  -0.16%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '683.00'
$ws.Range("E5").Value = '  -1.28%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '159.35'
$ws.Range("E6").Value = '  -1.87%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("E8").Value = '  -0.97%  '
$ws.Range("E9").Value = '  -1.15%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.16'
$ws.Range("E10").Value = '  -3.23%  '
$ws.Range("E11").Value = '  -0.13%  '
$ws.Range("E12").Value = '  -2.83%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.309.68'
$ws.Range("E13").Value = '  -0.18%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.47'
$ws.Range("E14").Value = '  -2.66%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.684.21'
$ws.Range("E15").Value = '  -0.23%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '69.326.70'
$ws.Range("E16").Value = '  -0.25%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.116'
$ws.Range("E17").Value = '  +1.97%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '16.09'
$ws.Range("E18").Value = '  -0.82%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.47'
$ws.Range("E19").Value = '  -1.57%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '469.13'
$ws.Range("E20").Value = '  -2.32%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.93'
$ws.Range("E21").Value = '  -0.53%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.654'
$ws.Range("E22").Value = '  -1.10%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '79.89'
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.834.08'
$ws.Range("E24").Value = '  -0.16%  '
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("E26").Value = '  -5.15%  '
$ws.Range("E27").Value = '  -3.59%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.16'
$ws.Range("E28").Value = '  -3.39%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.70'
$ws.Range("E29").Value = '  -0.96%  '
$ws.Range("E30").Value = '  -4.27%  '
$ws.Range("E31").Value = '  -2.76%  '
$ws.Range("E33").Value = '  -0.05%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.90'
$ws.Range("E34").Value = '  -0.30%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.675.76'
$ws.Range("E35").Value = '  +0.24%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.156'
$ws.Range("E36").Value = '  -6.40%  '
$ws.Range("E37").Value = '  -1.97%  '
$ws.Range("E38").Value = '  -0.51%  '
$ws.Range("E39").Value = '  +0.00%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.26'
$ws.Range("E40").Value = '  -2.84%  '
$ws.Range("E41").Value = '  -0.09%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0908'
$ws.Range("E42").Value = '  -2.25%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '169.85'
$ws.Range("E43").Value = '  +3.35%  '
$ws.Range("E44").Value = '  -1.20%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '47.67'
$ws.Range("E45").Value = '  -0.74%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '28.32'
$ws.Range("E46").Value = '  -6.17%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.12'
$ws.Range("E47").Value = '  -1.77%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.71'
$ws.Range("E48").Value = '  -3.26%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.30'
$ws.Range("E49").Value = '  -2.99%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.000277'
$ws.Range("E50").Value = '  -1.90%  '
$ws.Range("E51").Value = '  -3.39%  '
